$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.464.11'
Set-TextValue 'E2' '  -1.87%  '

Set-TextValue 'D3' '1.832.37'
Set-TextValue 'E3' '  -2.56%  '

Set-TextValue 'E4' '  -0.61%  '

Set-TextValue 'D5' '331.12'
Set-TextValue 'E5' '  -1.35%  '

Set-TextValue 'D6' '1.004'
Set-TextValue 'E6' '  -0.59%  '

Set-TextValue 'E7' '  -3.72%  '

Set-TextValue 'D8' '0.3816'
Set-TextValue 'E8' '  -3.37%  '

Set-TextValue 'D9' '46.72'
Set-TextValue 'E9' '  -0.90%  '

Set-TextValue 'E10' '  -1.38%  '

Set-TextValue 'D11' '0.9699'
Set-TextValue 'E11' '  -4.79%  '

Set-TextValue 'D12' '21.06'
Set-TextValue 'E12' '  -3.95%  '

Set-TextValue 'D13' '1.832.13'
Set-TextValue 'E13' '  -2.85%  '

Set-TextValue 'D14' '5.892'
Set-TextValue 'E14' '  -2.68%  '

Set-TextValue 'D15' '7.030'
Set-TextValue 'E15' '  -2.46%  '

Set-TextValue 'D16' '1.005'
Set-TextValue 'E16' '  -0.80%  '

Set-TextValue 'E17' '  -0.92%  '

Set-TextValue 'D18' '0.06649'
Set-TextValue 'E18' '  -1.13%  '

Set-TextValue 'E19' '  -2.09%  '

Set-TextValue 'D20' '17.01'
Set-TextValue 'E20' '  -0.41%  '

Set-TextValue 'D21' '1.004'
Set-TextValue 'E21' '  -0.52%  '

Set-TextValue 'D22' '27.462.69'
Set-TextValue 'E22' '  -1.84%  '

Set-TextValue 'D23' '5.343'
Set-TextValue 'E23' '  -3.13%  '

Set-TextValue 'D24' '10.81'
Set-TextValue 'E24' '  -1.86%  '

Set-TextValue 'E25' '  -1.45%  '

Set-TextValue 'D26' '2.052.40'
Set-TextValue 'E26' '  -2.60%  '

Set-TextValue 'D27' '157.49'
Set-TextValue 'E27' '  -0.65%  '

Set-TextValue 'E28' '  -2.66%  '

Set-TextValue 'E29' '  -2.33%  '

Set-TextValue 'D30' '5.302'
Set-TextValue 'E30' '  -3.61%  '

Set-TextValue 'D31' '119.03'
Set-TextValue 'E31' '  -2.14%  '

Set-TextValue 'E32' '  -2.78%  '

Set-TextValue 'D33' '0.09292'
Set-TextValue 'E33' '  -2.96%  '

Set-TextValue 'D34' '3.585'
Set-TextValue 'E34' '  -1.27%  '

Set-TextValue 'D35' '5.243'
Set-TextValue 'E35' '  -1.88%  '

Set-TextValue 'E36' '  -3.44%  '

Set-TextValue 'D37' '0.05927'
Set-TextValue 'E37' '  -2.53%  '

Set-TextValue 'D38' '0.02193'
Set-TextValue 'E38' '  -2.46%  '

Set-TextValue 'D39' '1.164'
Set-TextValue 'E39' '  -3.51%  '

Set-TextValue 'D40' '8.051'
Set-TextValue 'E40' '  -1.79%  '

Set-TextValue 'D41' '0.5790'
Set-TextValue 'E41' '  -3.53%  '

Set-TextValue 'E42' '  -3.28%  '

Set-TextValue 'E43' '  -3.21%  '

Set-TextValue 'D44' '1.248'
Set-TextValue 'E44' '  -0.98%  '

Set-TextValue 'D45' '0.5480'
Set-TextValue 'E45' '  -3.65%  '

Set-TextValue 'E46' '  -2.07%  '

Set-TextValue 'E47' '  -3.63%  '

Set-TextValue 'D48' '0.06641'
Set-TextValue 'E48' '  -2.41%  '

Set-TextValue 'D49' '110.33'
Set-TextValue 'E49' '  -2.18%  '

Set-TextValue 'E50' '  -3.12%  '

Set-TextValue 'D51' '1.003'
Set-TextValue 'E51' '  -0.67%  '
